$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths: only column C keeps a custom width; B & D go back to default ---
$ws.Columns("B").ColumnWidth = 8.43
$ws.Columns("D").ColumnWidth = 8.43
$ws.Columns("C").ColumnWidth = 14.59

# --- Header row (row 1) ---
$ws.Range("A1").Value = "user_id"
$ws.Range("B1").Value = "penjualan_kode"
$ws.Range("C1").Value = "pembeli"
$ws.Range("D1").Value = "barang_id"
$ws.Range("E1").Value = "jumlah"

$ws.Range("A1:E1").Font.Bold = $true
$ws.Range("A1:C1").HorizontalAlignment = -4108

# --- Data rows: clear previous date-format styling before writing new numbers ---
$ws.Range("D2:D3").Style = "Normal"
$ws.Range("E2:E3").Style = "Normal"

# Row 3's text is interned into the shared-string table before row 2's
# (and C3 before B3), matching the order in the target file.
$ws.Range("A3").Value = 3
$ws.Range("C3").Value = "Customer 16"
$ws.Range("B3").Value = "TXR0025"
$ws.Range("D3").Value = 10
$ws.Range("E3").Value = 1

$ws.Range("A2").Value = 3
$ws.Range("B2").Value = "TXR0024"
$ws.Range("C2").Value = "Customer 24"
$ws.Range("D2").Value = 10
$ws.Range("E2").Value = 1

# --- Row 4 becomes blank, but keep the cells present (stamped with a style) ---
$ws.Range("A4:E4").ClearContents()
$ws.Range("A4:E4").Style = "Normal"

# --- Selection matches the target ---
$ws.Range("D4:E4").Select()
